$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: values that are unambiguous text (contain 2+ dots) ---
# These can be assigned directly; Excel will not reinterpret them as numbers.
$ws.Range("D2").Value = "26.319.74"
$ws.Range("D3").Value = "1.678.07"
$ws.Range("D12").Value = "1.683.17"
$ws.Range("D15").Value = "1.908.99"
$ws.Range("D18").Value = "26.388.10"
$ws.Range("D39").Value = "1.096.18"
$ws.Range("D45").Value = "1.835.10"

# --- Column D: values that look like plain numbers ---
# Excel auto-converts numeric-looking text typed into a General cell into a
# real number, which would lose the original formatting (e.g. trailing zeros,
# "0.000008510"). To keep these as text (matching the original inlineStr cells)
# we stage the text in a scratch cell that is explicitly formatted as Text,
# then copy/paste-special only the *values* into the target cell so the target
# keeps its original (default) style while the content remains a text string.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "1.006"
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$scratch.Value = "218.93"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$scratch.Value = "0.5051"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$scratch.Value = "0.2660"
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$scratch.Value = "22.03"
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$scratch.Value = "0.06309"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$scratch.Value = "0.07368"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$scratch.Value = "4.521"
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$scratch.Value = "0.5779"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$scratch.Value = "0.000008510"
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$scratch.Value = "64.87"
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$scratch.Value = "5.000"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$scratch.Value = "10.86"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$scratch.Value = "186.32"
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$scratch.Value = "6.227"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$scratch.Value = "1.007"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$scratch.Value = "143.75"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$scratch.Value = "7.486"
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$scratch.Value = "0.1168"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$scratch.Value = "15.90"
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$scratch.Value = "1.346"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$scratch.Value = "0.05799"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$scratch.Value = "3.510"
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$scratch.Value = "3.502"
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$scratch.Value = "1.656"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Value = "1.008"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$scratch.Value = "0.5949"
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$scratch.Value = "2.362"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$scratch.Value = "2.675"
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$scratch.Value = "0.01601"
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$scratch.Value = "5.894"
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$scratch.Value = "0.8598"
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$scratch.Value = "0.00000000114"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$scratch.Value = "56.22"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$scratch.Value = "8.013"
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$scratch.Value = "0.4317"
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$scratch.Value = "0.05211"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$scratch.Clear()

# --- Column E: percentage-change text values ---
# These retain their leading/trailing double-space padding, which keeps Excel
# from treating them as numeric percentages, so a direct assignment is safe.
$ws.Range("E2").Value = "  -7.46%  "
$ws.Range("E3").Value = "  -6.21%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("E5").Value = "  -5.20%  "
$ws.Range("E6").Value = "  -13.84%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("E12").Value = "  -6.03%  "
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("E15").Value = "  -6.03%  "
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("E17").Value = "  -13.88%  "
$ws.Range("E18").Value = "  -7.18%  "
$ws.Range("E19").Value = "  -7.44%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  -4.66%  "
$ws.Range("E22").Value = "  -10.11%  "
$ws.Range("E23").Value = "  -8.03%  "
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("E25").Value = "  -5.49%  "
$ws.Range("E26").Value = "  -6.71%  "
$ws.Range("E27").Value = "  -6.87%  "
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("E30").Value = "  -5.33%  "
$ws.Range("E31").Value = "  -5.92%  "
$ws.Range("E32").Value = "  -6.75%  "
$ws.Range("E33").Value = "  -6.72%  "
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("E40").Value = "  -4.26%  "
$ws.Range("E41").Value = "  -6.39%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("E47").Value = "  -6.14%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("E50").Value = "  -3.42%  "
$ws.Range("E51").Value = "  -4.01%  "
